# DE OPSD 2020 smaller fixes and changes
# - extend the "plants" optional-attributes block (columns G:H) on sheet "raw"
#   with a new "availability" / "float64" row
# - add two new small reference tables on sheet "raw": "heatareas" (AK:AL)
#   and "demand_h" (AN:AO)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("raw")

# --- 1. plants: add trailing optional attribute "availability" (float64) ---
$ws.Range("G27").Value = "availability"
$ws.Range("H27").Value = "float64"

# --- 2. new "heatareas" table header (columns AK:AL) ------------------------
$ws.Range("AH1:AI1").Copy()
$ws.Range("AK1:AL1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AH2:AI2").Copy()
$ws.Range("AK2:AL2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AK1").Value = "heatareas"
$ws.Range("AK1:AL1").Merge()

$ws.Range("AK2").Value = "attributes "
$ws.Range("AL2").Value = "type"

$ws.Range("AK3").Value = "index"
$ws.Range("AL3").Value = "any"

# --- 3. new "demand_h" table (columns AN:AO) --------------------------------
$ws.Range("AH1:AI1").Copy()
$ws.Range("AN1:AO1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AH2:AI2").Copy()
$ws.Range("AN2:AO2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AH8:AI8").Copy()
$ws.Range("AN8:AO8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AN1").Value = "demand_h"
$ws.Range("AN1:AO1").Merge()

$ws.Range("AN2").Value = "attributes "
$ws.Range("AO2").Value = "type"

$ws.Range("AN3").Value = "index"
$ws.Range("AO3").Value = "any"

$ws.Range("AN4").Value = "timestep"
$ws.Range("AO4").Value = "any"

$ws.Range("AN5").Value = "heatarea"
$ws.Range("AO5").Value = "heatareas.index"

$ws.Range("AN6").Value = "demand_h"
$ws.Range("AO6").Value = "float64"

$ws.Range("AN8").Value = "optional attributes"

# --- 4. leave the sheet scrolled/selected where the edits were made --------
$ws.Activate()
$ws.Range("H1048576").Select()

Write-Host "applied heatareas/demand_h tables + plants.availability"
